# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210   (columns A..J)
#   *_new -> *_FV2304   (columns L..U)
# Then wrap the sheet's data range in a real Excel Table ("Table1") and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "old"/"new" header captions to the format-version captions ---
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1            # columns A..J
    $ws.Cells.Item(1, $col).Value = "$($oldHeaders[$i])_FV2210"
}

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 12            # columns L..U (K is "diff")
    $ws.Cells.Item(1, $col).Value = "$($oldHeaders[$i])_FV2304"
}

# --- 2. Turn the used range into an Excel Table so the new headers become ---
#        proper table column names (ListObject / Table1).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U88"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (View > Freeze Panes > Freeze Top Row) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
